# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
# (commit: "Updated cryptos list on Sun May 21 19:37:05 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '27.094.23'; ForceText = $false },
    @{ Cell = 'E2'; Value = '  -1.21%  '; ForceText = $false },
    @{ Cell = 'D3'; Value = '1.822.79'; ForceText = $false },
    @{ Cell = 'E3'; Value = '  -1.21%  '; ForceText = $false },
    @{ Cell = 'E4'; Value = '  -0.61%  '; ForceText = $false },
    @{ Cell = 'D5'; Value = '311.55'; ForceText = $true },
    @{ Cell = 'E5'; Value = '  -1.57%  '; ForceText = $false },
    @{ Cell = 'D6'; Value = '1.008'; ForceText = $true },
    @{ Cell = 'E6'; Value = '  -0.55%  '; ForceText = $false },
    @{ Cell = 'E7'; Value = '  -2.34%  '; ForceText = $false },
    @{ Cell = 'D8'; Value = '0.3640'; ForceText = $true },
    @{ Cell = 'E8'; Value = '  -1.71%  '; ForceText = $false },
    @{ Cell = 'D9'; Value = '0.07289'; ForceText = $true },
    @{ Cell = 'D10'; Value = '0.8700'; ForceText = $true },
    @{ Cell = 'E10'; Value = '  -1.98%  '; ForceText = $false },
    @{ Cell = 'D11'; Value = '20.14'; ForceText = $true },
    @{ Cell = 'E11'; Value = '  -1.88%  '; ForceText = $false },
    @{ Cell = 'D12'; Value = '1.862.89'; ForceText = $false },
    @{ Cell = 'E12'; Value = '  -0.33%  '; ForceText = $false },
    @{ Cell = 'D13'; Value = '0.07582'; ForceText = $true },
    @{ Cell = 'E13'; Value = '  +2.53%  '; ForceText = $false },
    @{ Cell = 'D14'; Value = '5.349'; ForceText = $true },
    @{ Cell = 'E14'; Value = '  -2.62%  '; ForceText = $false },
    @{ Cell = 'D15'; Value = '92.56'; ForceText = $true },
    @{ Cell = 'E15'; Value = '  -0.77%  '; ForceText = $false },
    @{ Cell = 'D16'; Value = '6.473'; ForceText = $true },
    @{ Cell = 'E16'; Value = '  -1.81%  '; ForceText = $false },
    @{ Cell = 'D17'; Value = '1.009'; ForceText = $true },
    @{ Cell = 'E17'; Value = '  -0.57%  '; ForceText = $false },
    @{ Cell = 'D18'; Value = '0.000008640'; ForceText = $true },
    @{ Cell = 'E18'; Value = '  -2.55%  '; ForceText = $false },
    @{ Cell = 'D20'; Value = '27.369.98'; ForceText = $false },
    @{ Cell = 'E20'; Value = '  -0.24%  '; ForceText = $false },
    @{ Cell = 'E21'; Value = '  -2.64%  '; ForceText = $false },
    @{ Cell = 'D22'; Value = '5.202'; ForceText = $true },
    @{ Cell = 'E22'; Value = '  -2.66%  '; ForceText = $false },
    @{ Cell = 'E23'; Value = '  -1.49%  '; ForceText = $false },
    @{ Cell = 'D24'; Value = '2.093.42'; ForceText = $false },
    @{ Cell = 'E24'; Value = '  +0.81%  '; ForceText = $false },
    @{ Cell = 'B25'; Value = 'Monero'; ForceText = $false },
    @{ Cell = 'C25'; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; ForceText = $false },
    @{ Cell = 'D25'; Value = '151.65'; ForceText = $true },
    @{ Cell = 'E25'; Value = '  -0.60%  '; ForceText = $false },
    @{ Cell = 'B26'; Value = 'Toncoin'; ForceText = $false },
    @{ Cell = 'C26'; Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; ForceText = $false },
    @{ Cell = 'D26'; Value = '1.872'; ForceText = $true },
    @{ Cell = 'E26'; Value = '  -2.06%  '; ForceText = $false },
    @{ Cell = 'E27'; Value = '  -2.40%  '; ForceText = $false },
    @{ Cell = 'D28'; Value = '2.101'; ForceText = $true },
    @{ Cell = 'E28'; Value = '  -3.49%  '; ForceText = $false },
    @{ Cell = 'D29'; Value = '116.04'; ForceText = $true },
    @{ Cell = 'E29'; Value = '  -1.77%  '; ForceText = $false },
    @{ Cell = 'D30'; Value = '5.066'; ForceText = $true },
    @{ Cell = 'E30'; Value = '  -4.26%  '; ForceText = $false },
    @{ Cell = 'D31'; Value = '0.08914'; ForceText = $true },
    @{ Cell = 'E31'; Value = '  -0.75%  '; ForceText = $false },
    @{ Cell = 'E32'; Value = '  +0.41%  '; ForceText = $false },
    @{ Cell = 'D33'; Value = '0.7339'; ForceText = $true },
    @{ Cell = 'E33'; Value = '  -3.70%  '; ForceText = $false },
    @{ Cell = 'D34'; Value = '4.450'; ForceText = $true },
    @{ Cell = 'E34'; Value = '  -2.65%  '; ForceText = $false },
    @{ Cell = 'E35'; Value = '  -3.42%  '; ForceText = $false },
    @{ Cell = 'D36'; Value = '1.009'; ForceText = $true },
    @{ Cell = 'E36'; Value = '  -0.49%  '; ForceText = $false },
    @{ Cell = 'D37'; Value = '2.528'; ForceText = $true },
    @{ Cell = 'E37'; Value = '  +5.48%  '; ForceText = $false },
    @{ Cell = 'E38'; Value = '  -3.11%  '; ForceText = $false },
    @{ Cell = 'D39'; Value = '0.05256'; ForceText = $true },
    @{ Cell = 'D40'; Value = '0.01921'; ForceText = $true },
    @{ Cell = 'E40'; Value = '  -2.42%  '; ForceText = $false },
    @{ Cell = 'D41'; Value = '2.934'; ForceText = $true },
    @{ Cell = 'E41'; Value = '  -2.54%  '; ForceText = $false },
    @{ Cell = 'D42'; Value = '7.155'; ForceText = $true },
    @{ Cell = 'E42'; Value = '  -2.30%  '; ForceText = $false },
    @{ Cell = 'D43'; Value = '0.5212'; ForceText = $true },
    @{ Cell = 'E43'; Value = '  -2.64%  '; ForceText = $false },
    @{ Cell = 'E44'; Value = '  -2.14%  '; ForceText = $false },
    @{ Cell = 'D45'; Value = '8.271'; ForceText = $true },
    @{ Cell = 'E45'; Value = '  -3.19%  '; ForceText = $false },
    @{ Cell = 'D46'; Value = '0.4890'; ForceText = $true },
    @{ Cell = 'E46'; Value = '  -1.51%  '; ForceText = $false },
    @{ Cell = 'E47'; Value = '  -0.63%  '; ForceText = $false },
    @{ Cell = 'B48'; Value = 'Quant'; ForceText = $false },
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; ForceText = $false },
    @{ Cell = 'D48'; Value = '104.11'; ForceText = $true },
    @{ Cell = 'E48'; Value = '  -0.92%  '; ForceText = $false },
    @{ Cell = 'B49'; Value = 'EnergySwap'; ForceText = $false },
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; ForceText = $false },
    @{ Cell = 'D49'; Value = '10.14'; ForceText = $true },
    @{ Cell = 'E49'; Value = '  -4.04%  '; ForceText = $false },
    @{ Cell = 'D50'; Value = '1.635'; ForceText = $true },
    @{ Cell = 'E50'; Value = '  -2.96%  '; ForceText = $false },
    @{ Cell = 'D51'; Value = '0.06254'; ForceText = $true },
    @{ Cell = 'E51'; Value = '  -1.28%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Value looks like a plain number ("311.55", "0.3640", ...) - the sheet
        # stores these as text (inline strings), so force Text format, assign,
        # then clear the format again so the cell style index is unaffected.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.ClearFormats()
    } else {
        $rng.Value = $u.Value
    }
}

